$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$fileData = New-Object 'object[,]' 185,1
$dateData = New-Object 'object[,]' 185,1
$fileData[0,0] = 'Routine_Care/Nursing for Arterial and Central Venous Lines.pdf'
$dateData[0,0] = '2011-03-19'
$fileData[1,0] = 'Routine_Care/VTE_Prevention/TED Stocking Sizing.pdf'
$dateData[1,0] = '2011-06-19'
$fileData[2,0] = 'Routine_Care/Faecal  incontinence skin care.pdf'
$dateData[2,0] = '2011-07-19'
$fileData[3,0] = 'Breathing(Respiratory)/Equipment/IPPB using an ICU Ventilator.pdf'
$dateData[3,0] = '2011-12-19'
$fileData[4,0] = 'Drugs/heparin_critical_care_only.pdf'
$dateData[4,0] = '2014-09-19'
$fileData[5,0] = 'Routine_Care/Invasive Flush Systems.pdf'
$dateData[5,0] = '2014-10-19'
$fileData[6,0] = 'GI_Liver_and_Transplant/Pancreatic Irrigation.pdf'
$dateData[6,0] = '2014-11-19'
$fileData[7,0] = 'Neurological/Thiopentone levels.pdf'
$dateData[7,0] = '2014-12-19'
$fileData[8,0] = 'Infection_and_sepsis/Ebola/Ebola.pdf'
$dateData[8,0] = '2015-03-19'
$fileData[9,0] = 'GI_Liver_and_Transplant/Nasal bridle.pdf'
$dateData[9,0] = '2016-01-19'
$fileData[10,0] = 'Breathing(Respiratory)/Equipment/AMBU AScope.pdf'
$dateData[10,0] = '2016-04-19'
$fileData[11,0] = 'Cardiovascular/EZ-IO Intraosseus Access Device_pub_em.pdf'
$dateData[11,0] = '2016-08-19'
$fileData[12,0] = 'Routine_Care/Central venous catheter removal.pdf'
$dateData[12,0] = '2017-01-19'
$fileData[13,0] = 'Routine_Care/Tracheostomy_nursing_care.pdf'
$dateData[13,0] = '2017-08-19'
$fileData[14,0] = 'Airway/Critical care extubation checklist.pdf'
$dateData[14,0] = '2017-10-19'
$fileData[15,0] = 'ECLS/Extra Corporeal Carbon Dioxide Removal.pdf'
$dateData[15,0] = '2017-10-19'
$fileData[16,0] = 'Airway/Tracheostomy_Laryngectomy/Hospital_in-patients_with_a_Tracheostomy.pdf'
$dateData[16,0] = '2017-11-19'
$fileData[17,0] = 'Trauma and Burns/Mangement of burns.pdf'
$dateData[17,0] = '2018-05-19'
$fileData[18,0] = 'Drugs/diazepam_diazemuls.pdf'
$dateData[18,0] = '2018-07-19'
$fileData[19,0] = 'End_of_life_care/Reasons to report a death to PF.pdf'
$dateData[19,0] = '2019-05-19'
$fileData[20,0] = 'Breathing(Respiratory)/salbutamol and ipratroprium MDI.pdf'
$dateData[20,0] = '2019-05-19'
$fileData[21,0] = 'Airway/Tracheostomy_Laryngectomy/Tracheostomy suctioning cleaning guideline.pdf'
$dateData[21,0] = '2019-06-19'
$fileData[22,0] = 'Neurological/SOP -  Femoral site care.pdf'
$dateData[22,0] = '2019-06-19'
$fileData[23,0] = 'Airway/Tracheostomy_Laryngectomy/Tracheostomy change in Critical Care.pdf'
$dateData[23,0] = '2019-06-19'
$fileData[24,0] = 'Drugs/sodium_bicarbonate.pdf'
$dateData[24,0] = '2019-08-19'
$fileData[25,0] = 'Post_op_care/Anticoagulation antiplatelet agents and epidural analgesia.pdf'
$dateData[25,0] = '2019-09-19'
$fileData[26,0] = 'Post_op_care/Epidural top-up.pdf'
$dateData[26,0] = '2020-01-19'
$fileData[27,0] = 'Drugs/heparin for Haemofiltration.pdf'
$dateData[27,0] = '2020-03-19'
$fileData[28,0] = 'Covid-19/SJH/SJH COVID19 ED Intubation Action Card.pdf'
$dateData[28,0] = '2020-03-19'
$fileData[29,0] = 'Covid-19/SJH/SJH COVID19 ITU Intubation Action Card.pdf'
$dateData[29,0] = '2020-03-19'
$fileData[30,0] = 'Covid-19/WGH/CoVid intubation checklist WGH.pdf'
$dateData[30,0] = '2020-03-19'
$fileData[31,0] = 'Airway/Emergency intubation checklist_em_pub.pdf'
$dateData[31,0] = '2020-03-19'
$fileData[32,0] = 'Diabetes_and_Glucose/Hyperosmolar Hyperglycaemic State.pdf'
$dateData[32,0] = '2020-03-19'
$fileData[33,0] = 'Drugs/fentanyl.pdf'
$dateData[33,0] = '2020-04-19'
$fileData[34,0] = 'Airway/Tracheostomy_Laryngectomy/Tracheostomy guideline.pdf'
$dateData[34,0] = '2020-05-19'
$fileData[35,0] = 'Covid-19/WGH/WGH_CT_Transfer_May.pdf'
$dateData[35,0] = '2020-07-19'
$fileData[36,0] = 'Cardiovascular/GJNH Acute Heart Failure Referral Form.pdf'
$dateData[36,0] = '2020-08-19'
$fileData[37,0] = 'Organ_donation/Donation after circulatory death.pdf'
$dateData[37,0] = '2020-11-19'
$fileData[38,0] = 'Airway/Percutaneous tracheostomy checklist.pdf'
$dateData[38,0] = '2021-02-19'
$fileData[39,0] = 'Delirium/Managing a Potentially Violent Patient.pdf'
$dateData[39,0] = '2021-05-19'
$fileData[40,0] = 'Delirium/Risk assessment posi mit.pdf'
$dateData[40,0] = '2021-05-19'
$fileData[41,0] = 'Infection_and_sepsis/SOP Ultrasound Cleaning.pdf'
$dateData[41,0] = '2021-05-19'
$fileData[42,0] = 'Breathing(Respiratory)/HFNO.pdf'
$dateData[42,0] = '2021-06-19'
$fileData[43,0] = 'Delirium/Drugs Causing Delirium and Agitiation.pdf'
$dateData[43,0] = '2021-06-19'
$fileData[44,0] = 'Drugs/ketamine_in_asthma.pdf'
$dateData[44,0] = '2021-06-19'
$fileData[45,0] = 'Airway/McGrath Mac.pdf'
$dateData[45,0] = '2021-06-19'
$fileData[46,0] = 'Airway/Tracheostomy_Laryngectomy/Tracheostomy safety box contents.pdf'
$dateData[46,0] = '2021-06-19'
$fileData[47,0] = 'GI_Liver_and_Transplant/Treatment of constipation.pdf'
$dateData[47,0] = '2021-06-19'
$fileData[48,0] = 'GI_Liver_and_Transplant/Abdominal pressure measurement.pdf'
$dateData[48,0] = '2021-06-19'
$fileData[49,0] = 'Neurological/Sub arachnoid haemorrhage management.pdf'
$dateData[49,0] = '2021-06-19'
$fileData[50,0] = 'Airway/Anticipated difficult airway tool.pdf'
$dateData[50,0] = '2021-06-19'
$fileData[51,0] = 'End_of_life_care/Documentation following death.pdf'
$dateData[51,0] = '2021-09-19'
$fileData[52,0] = 'Drugs/zanamivir.pdf'
$dateData[52,0] = '2021-12-19'
$fileData[53,0] = 'Routine_Care/bBraun Spaceplus Failure EMERGENCY ACTION CARD_em.pdf'
$dateData[53,0] = '2022-01-19'
$fileData[54,0] = 'Breathing(Respiratory)/Equipment/HFNO Set Up.pdf'
$dateData[54,0] = '2022-03-19'
$fileData[55,0] = 'Drugs/insulin.pdf'
$dateData[55,0] = '2022-03-19'
$fileData[56,0] = 'Breathing(Respiratory)/Inhaled Nitrous Oxide.pdf'
$dateData[56,0] = '2022-04-19'
$fileData[57,0] = 'Cardiovascular/Steroids for Septic Shock.pdf'
$dateData[57,0] = '2022-05-19'
$fileData[58,0] = 'Breathing(Respiratory)/Equipment/APRV.pdf'
$dateData[58,0] = '2022-05-19'
$fileData[59,0] = 'Breathing(Respiratory)/Equipment/T piece Y piece.pdf'
$dateData[59,0] = '2022-06-19'
$fileData[60,0] = 'Post_op_care/Epidural Haematoma.pdf'
$dateData[60,0] = '2022-06-19'
$fileData[61,0] = 'Neurological/SOP for review of Neurosurgical patients in ITU by neurosurgical team.pdf'
$dateData[61,0] = '2022-06-19'
$fileData[62,0] = 'Policies_and_admin/General Critical Care Interaction with HEPMA_pub.pdf'
$dateData[62,0] = '2022-07-19'
$fileData[63,0] = 'Drugs/midazolam and thiopental levels.pdf'
$dateData[63,0] = '2022-08-19'
$fileData[64,0] = 'Breathing(Respiratory)/Equipment/HFNO through ventilator.pdf'
$dateData[64,0] = '2022-10-19'
$fileData[65,0] = 'Routine_Care/VTE_Prevention/Dalteparin_thromboprophylaxis.pdf'
$dateData[65,0] = '2022-11-19'
$fileData[66,0] = 'Post_op_care/Adult Scoliosis Spinal Surgery Post-Op Care.pdf'
$dateData[66,0] = '2022-11-19'
$fileData[67,0] = 'Post_op_care/Post op care pharyngo-laryngo-oesphagectomy PLOG.pdf'
$dateData[67,0] = '2022-12-19'
$fileData[68,0] = 'GI_Liver_and_Transplant/Nasogastric feeding protocol.pdf'
$dateData[68,0] = '2023-01-19'
$fileData[69,0] = 'Diabetes_and_Glucose/Intravenous Insulin Therapy (not for DKA or HHS).pdf'
$dateData[69,0] = '2023-02-19'
$fileData[70,0] = 'Drugs/Antibiotic doses in CVVHD.pdf'
$dateData[70,0] = '2023-02-19'
$fileData[71,0] = 'ECLS/RIE ECLS Anti Xa Protocol.pdf'
$dateData[71,0] = '2023-04-19'
$fileData[72,0] = 'GI_Liver_and_Transplant/Jejunostomy feeding protocol.pdf'
$dateData[72,0] = '2023-04-19'
$fileData[73,0] = 'GI_Liver_and_Transplant/Nasojejunal feeding protocol.pdf'
$dateData[73,0] = '2023-04-19'
$fileData[74,0] = 'Infection_and_sepsis/Winter Infections Stepdown Guidance.pdf'
$dateData[74,0] = '2023-05-19'
$fileData[75,0] = 'Drugs/vasopressin organ donation.pdf'
$dateData[75,0] = '2023-05-19'
$fileData[76,0] = 'Drugs/vasopressin_sepsis.pdf'
$dateData[76,0] = '2023-05-19'
$fileData[77,0] = 'Covid-19/videos/Donning and Doffing Video.pdf'
$dateData[77,0] = '2023-06-19'
$fileData[78,0] = 'Transfer/ACCP Transfers.pdf'
$dateData[78,0] = '2023-06-19'
$fileData[79,0] = 'Breathing(Respiratory)/CPAP.pdf'
$dateData[79,0] = '2023-07-19'
$fileData[80,0] = 'Breathing(Respiratory)/Equipment/Ventilators Circuits Filters and Closed Suction - Set up and Maintenance.pdf'
$dateData[80,0] = '2023-07-19'
$fileData[81,0] = 'Infection_and_sepsis/Infection indications for IVIG.pdf'
$dateData[81,0] = '2023-07-19'
$fileData[82,0] = 'Drugs/piperacillin_tazobactam extended_infusion.pdf'
$dateData[82,0] = '2023-07-19'
$fileData[83,0] = 'Breathing(Respiratory)/Equipment/Bipap V60.pdf'
$dateData[83,0] = '2023-07-19'
$fileData[84,0] = 'Covid-19/Covid 19 Death Certification Guideline.pdf'
$dateData[84,0] = '2023-08-19'
$fileData[85,0] = 'Procedures/CVC Guidance/Securing CVCs.pdf'
$dateData[85,0] = '2023-08-19'
$fileData[86,0] = 'Routine_Care/Video Communication.pdf'
$dateData[86,0] = '2023-09-19'
$fileData[87,0] = 'Neurological/Treatment of status epilepticus.pdf'
$dateData[87,0] = '2023-09-19'
$fileData[88,0] = 'Drugs/isoprenaline.pdf'
$dateData[88,0] = '2023-10-19'
$fileData[89,0] = 'Cardiovascular/Cardiogenic Shock.pdf'
$dateData[89,0] = '2023-10-19'
$fileData[90,0] = 'Haematology_CAR-T/Haem_ICU_transfer.pdf'
$dateData[90,0] = '2024-01-19'
$fileData[91,0] = 'Cardiovascular/Management of hypertension within Critical Care.pdf'
$dateData[91,0] = '2024-02-19'
$fileData[92,0] = 'Drugs/aminophylline.pdf'
$dateData[92,0] = '2024-02-19'
$fileData[93,0] = 'Haematology_CAR-T/ICANS.pdf'
$dateData[93,0] = '2024-03-19'
$fileData[94,0] = 'Drugs/rocuronium.pdf'
$dateData[94,0] = '2024-03-19'
$fileData[95,0] = 'Drugs/pancuronium.pdf'
$dateData[95,0] = '2024-03-19'
$fileData[96,0] = 'Drugs/phenytoin.pdf'
$dateData[96,0] = '2024-03-19'
$fileData[97,0] = 'Haematology_CAR-T/CRS.pdf'
$dateData[97,0] = '2024-03-19'
$fileData[98,0] = 'Drugs/milrinone.pdf'
$dateData[98,0] = '2024-04-19'
$fileData[99,0] = 'Policies_and_admin/General Critical Care SOP_pub.pdf'
$dateData[99,0] = '2024-04-19'
$fileData[100,0] = 'Covid-19/COVID 19 ICM guidance basic goals_June_2022.pdf'
$dateData[100,0] = '2024-05-19'
$fileData[101,0] = 'Neurological/Critical Care MRI Procedure_pub.pdf'
$dateData[101,0] = '2024-05-19'
$fileData[102,0] = 'Infection_and_sepsis/Trip Out of Unit infection guidance.pdf'
$dateData[102,0] = '2024-05-19'
$fileData[103,0] = 'Organ_donation/Organ Retrieval SOP.pdf'
$dateData[103,0] = '2024-05-19'
$fileData[104,0] = 'Drugs/clonidine.pdf'
$dateData[104,0] = '2024-05-19'
$fileData[105,0] = 'Neurological/Management of traumatic brain injury.pdf'
$dateData[105,0] = '2024-05-19'
$fileData[106,0] = 'Ethics_and_Law/Care at the End of Life (FICM).pdf'
$dateData[106,0] = '2024-05-19'
$fileData[107,0] = 'Ethics_and_Law/DNACPR policy for Scotland.pdf'
$dateData[107,0] = '2024-05-19'
$fileData[108,0] = 'GI_Liver_and_Transplant/ICU - Upper GI bleeding (Endoscopy guideline).pdf'
$dateData[108,0] = '2024-05-19'
$fileData[109,0] = 'End_of_life_care/CMO & NRS Guidance for Doctors completing MCCD.pdf'
$dateData[109,0] = '2024-05-19'
$fileData[110,0] = 'End_of_life_care/Palliative extubation & withdrawal of invasive ventilatory support nursing checklist.pdf'
$dateData[110,0] = '2024-05-19'
$fileData[111,0] = 'Airway/Cook Staged Extubation Set.pdf'
$dateData[111,0] = '2024-06-02'
$fileData[112,0] = 'Drugs/noradrenaline (central).pdf'
$dateData[112,0] = '2024-06-19'
$fileData[113,0] = 'Post_op_care/Epidural hypotension.pdf'
$dateData[113,0] = '2024-06-26'
$fileData[114,0] = 'Breathing(Respiratory)/Equipment/Passy Muir Valve.pdf'
$dateData[114,0] = '2024-07-19'
$fileData[115,0] = 'GI_Liver_and_Transplant/Fulminant Liver Failure.pdf'
$dateData[115,0] = '2024-07-19'
$fileData[116,0] = 'Drugs/dexmedetomidine.pdf'
$dateData[116,0] = '2024-07-19'
$fileData[117,0] = 'Drugs/glyceryl_trinitrate.pdf'
$dateData[117,0] = '2024-07-19'
$fileData[118,0] = 'GI_Liver_and_Transplant/Confirmation of Nasogastric Tube Position.pdf'
$dateData[118,0] = '2024-07-19'
$fileData[119,0] = 'Infection_and_sepsis/Antifungal guidance in critical care.pdf'
$dateData[119,0] = '2024-07-25'
$fileData[120,0] = 'Breathing(Respiratory)/ARDS Strategy.pdf'
$dateData[120,0] = '2024-08-15'
$fileData[121,0] = 'Transfer/Transfer Outdoors to Garden Guideline.pdf'
$dateData[121,0] = '2024-08-19'
$fileData[122,0] = 'Cardiovascular/Intra Aortic Balloon Pump Bedside Checks_pub.pdf'
$dateData[122,0] = '2024-08-19'
$fileData[123,0] = 'Cardiovascular/Intra Aortic Balloon Pump Guideline_pub.pdf'
$dateData[123,0] = '2024-08-19'
$fileData[124,0] = 'Drugs/dobutamine.pdf'
$dateData[124,0] = '2024-10-19'
$fileData[125,0] = 'Drugs/adrenaline.pdf'
$dateData[125,0] = '2024-10-19'
$fileData[126,0] = 'Drugs/hydralazine.pdf'
$dateData[126,0] = '2024-10-24'
$fileData[127,0] = 'Post_op_care/Major OMFS Free Flap.pdf'
$dateData[127,0] = '2024-11-20'
$fileData[128,0] = 'Drugs/Alteplase for massive PE.pdf'
$dateData[128,0] = '2024-11-24'
$fileData[129,0] = 'Drugs/alfentanil.pdf'
$dateData[129,0] = '2024-11-24'
$fileData[130,0] = 'Drugs/magnesium.pdf'
$dateData[130,0] = '2024-12-24'
$fileData[131,0] = 'Drugs/ICU - IV drug infusions.pdf'
$dateData[131,0] = '2025-01-19'
$fileData[132,0] = 'Drugs/neostigmine.pdf'
$dateData[132,0] = '2025-01-19'
$fileData[133,0] = 'Drugs/vancomycin.pdf'
$dateData[133,0] = '2025-01-19'
$fileData[134,0] = 'Drugs/labetalol.pdf'
$dateData[134,0] = '2025-02-19'
$fileData[135,0] = 'Neurological/Intrathecal policy RIE.pdf'
$dateData[135,0] = '2025-02-19'
$fileData[136,0] = 'Infection_and_sepsis/Initial investigation and management in unidentified Infections.pdf'
$dateData[136,0] = '2025-02-19'
$fileData[137,0] = 'Drugs/midazolam.pdf'
$dateData[137,0] = '2025-03-19'
$fileData[138,0] = 'Cardiovascular/Management of Acute Type B Aortic Dissection Guideline.pdf'
$dateData[138,0] = '2025-03-19'
$fileData[139,0] = 'Drugs/potassium.pdf'
$dateData[139,0] = '2025-04-19'
$fileData[140,0] = 'Drugs/nimodipine.pdf'
$dateData[140,0] = '2025-04-19'
$fileData[141,0] = 'Drugs/salbutamol.pdf'
$dateData[141,0] = '2025-04-19'
$fileData[142,0] = 'Procedures/CVC Guidance/CVC NHL  April 2023.pdf'
$dateData[142,0] = '2025-04-19'
$fileData[143,0] = 'Routine_Care/ICU Eye Care Guideline.pdf'
$dateData[143,0] = '2025-05-19'
$fileData[144,0] = 'Drugs/nicardipine.pdf'
$dateData[144,0] = '2025-05-19'
$fileData[145,0] = 'Drugs/phenobarbitone.pdf'
$dateData[145,0] = '2025-05-19'
$fileData[146,0] = 'Procedures/Arterial Line insertion for ACCPs.pdf'
$dateData[146,0] = '2025-05-19'
$fileData[147,0] = 'Drugs/amiodarone.pdf'
$dateData[147,0] = '2025-05-19'
$fileData[148,0] = 'Drugs/morphine.pdf'
$dateData[148,0] = '2025-06-19'
$fileData[149,0] = 'Drugs/phenylephrine.pdf'
$dateData[149,0] = '2025-06-19'
$fileData[150,0] = 'Drugs/noradrenaline (peripheral).pdf'
$dateData[150,0] = '2025-06-19'
$fileData[151,0] = 'Breathing(Respiratory)/Manual Ventilation and MHI.pdf'
$dateData[151,0] = '2025-06-19'
$fileData[152,0] = 'Drugs/Epoprostenol.pdf'
$dateData[152,0] = '2025-06-19'
$fileData[153,0] = 'Neurological/Ventriculitis Guideline.pdf'
$dateData[153,0] = '2025-06-19'
$fileData[154,0] = 'Drugs/calcium.pdf'
$dateData[154,0] = '2025-07-19'
$fileData[155,0] = 'Cardiovascular/Cardiac Output Monitoring_pub .pdf'
$dateData[155,0] = '2025-07-19'
$fileData[156,0] = 'Cardiovascular/Pulmonary_Embolism_and_DVT/Catheter directed thrombolysis of iliofemoral DVT alteplase_pub.pdf'
$dateData[156,0] = '2025-07-19'
$fileData[157,0] = 'Drugs/atracurium.pdf'
$dateData[157,0] = '2025-08-19'
$fileData[158,0] = 'Airway/Tracheostomy_Laryngectomy/Decannulation Guidline.pdf'
$dateData[158,0] = '2025-08-19'
$fileData[159,0] = 'Drugs/Vancomycin Continuous Infusion Fluid Restricted.pdf'
$dateData[159,0] = '2025-08-19'
$fileData[160,0] = 'Drugs/dalteparin_thromboprophylaxis.pdf'
$dateData[160,0] = '2025-08-19'
$fileData[161,0] = 'Drugs/valproate.pdf'
$dateData[161,0] = '2025-10-19'
$fileData[162,0] = 'Policies_and_admin/Pet Visitation.pdf'
$dateData[162,0] = '2025-10-19'
$fileData[163,0] = 'Procedures/Inadvertent Catheter Placement Guideline.pdf'
$dateData[163,0] = '2025-10-19'
$fileData[164,0] = 'Drugs/stress ulcer prophylaxis.pdf'
$dateData[164,0] = '2025-11-19'
$fileData[165,0] = 'Drugs/Phosphate.pdf'
$dateData[165,0] = '2025-11-19'
$fileData[166,0] = 'Drugs/ketamine_for_status epilepticus.pdf'
$dateData[166,0] = '2025-11-19'
$fileData[167,0] = 'GI_Liver_and_Transplant/Prokinetics in ICU.pdf'
$dateData[167,0] = '2026-01-19'
$fileData[168,0] = 'Drugs/Octreotide.pdf'
$dateData[168,0] = '2026-01-19'
$fileData[169,0] = 'Breathing(Respiratory)/Proning Guideline.pdf'
$dateData[169,0] = '2026-01-19'
$fileData[170,0] = 'Drugs/Thiopentone.pdf'
$dateData[170,0] = '2026-01-19'
$fileData[171,0] = 'Procedures/ACCPs acquiring initial CVC competencies.pdf'
$dateData[171,0] = '2026-03-19'
$fileData[172,0] = 'Post_op_care/Prevention and treatment of paraplegia after major aortic procedures.pdf'
$dateData[172,0] = '2026-03-19'
$fileData[173,0] = 'Procedures/ACCP CVC placement following completion of initial competencies.pdf'
$dateData[173,0] = '2026-03-19'
$fileData[174,0] = 'Delirium/Violence and Agression.pdf'
$dateData[174,0] = '2026-05-19'
$fileData[175,0] = 'Post_op_care/Care of the Transgender Patient.pdf'
$dateData[175,0] = '2026-07-19'
$fileData[176,0] = 'GI_Liver_and_Transplant/Plasma exchange in Acute Liver Failure.pdf'
$dateData[176,0] = '2026-11-19'
$fileData[177,0] = 'End_of_life_care/Guideline following Sudden Cardiac Death where death occurs in ICU.pdf'
$dateData[177,0] = '2027-01-19'
$fileData[178,0] = 'Breathing(Respiratory)/Equipment/NIV through Drager Vent Set up in Critical Care.pdf'
$dateData[178,0] = '2027-01-19'
$fileData[179,0] = 'Transfer/Transfer Guidelines.pdf'
$dateData[179,0] = '2027-02-19'
$fileData[180,0] = 'Breathing(Respiratory)/Equipment/NIV through Nihon Kohden  Setup.pdf'
$dateData[180,0] = '2028-02-19'
$fileData[181,0] = 'Infection_and_sepsis/Influenza in Critical Care.pdf'
$dateData[181,0] = '2028-05-19'
$fileData[182,0] = 'Policies_and_admin/Anticipated Post op flow surgical patients.pdf'
$dateData[182,0] = '2028-07-19'
$fileData[183,0] = 'Policies_and_admin/Discharge Home from Critical Care.pdf'
$dateData[183,0] = '2028-11-19'
$fileData[184,0] = 'Policies_and_admin/Repatriaiton Checklist for Critical Care.pdf'
$dateData[184,0] = '2030-11-19'

$ws.Range("A2:A186").Value = $fileData

$dateRange = $ws.Range("B2:B186")
$dateRange.NumberFormat = "@"
$dateRange.Value = $dateData
$dateRange.ClearFormats()
